$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "Лаб1" in D3
$ws.Range("D3").Value = "Лаб1"

# Set grade 5 for the students that received the new "Лаб1" grade
$rows = @(5,6,8,10,13,17,23,24,26,27,28,30)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = 5
}

# Update the frozen-pane top-left cell and active selection to reflect scrolled view
$ws.Application.ActiveWindow.ScrollRow = 18
$ws.Range("D30").Select()
